$d = $word.ActiveDocument

# Item 1: '1. Download postgres database' -> split out 'postgres' with spell-check markers
$searchRange = $d.Content
$found = $searchRange.Find.Execute('1. Download postgres database', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find marker: 1. Download postgres database" }
$targetPara = $searchRange.Paragraphs(1)
$targetPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">1. Download </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>postgres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> database</w:t></w:r><w:r><w:t xml:space="preserve"> and instal</w:t></w:r><w:r><w:t>l</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Item 2: 'Download dbeaver.' -> split out 'dbeaver' with spell-check markers
$searchRange = $d.Content
$found = $searchRange.Find.Execute('Download dbeaver.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find marker: Download dbeaver." }
$targetPara = $searchRange.Paragraphs(1)
$targetPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">2. </w:t></w:r><w:r><w:t xml:space="preserve">Download </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dbeaver</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Item 3: 'Launch dbeaver and create new postgres database instance' -> split dbeaver/postgres
$searchRange = $d.Content
$found = $searchRange.Find.Execute('Launch dbeaver and c', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find marker: Launch dbeaver and c" }
$targetPara = $searchRange.Paragraphs(1)
$targetPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">3. </w:t></w:r><w:r><w:t xml:space="preserve">Launch </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dbeaver</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and c</w:t></w:r><w:r><w:t xml:space="preserve">reate new </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>postgres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> database instance</w:t></w:r><w:r><w:t xml:space="preserve"> and schema</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Code block: wrap jdbc/postgres/username/password terms in proofErr markers
$searchRange = $d.Content
$found = $searchRange.Find.Execute('spring.datasource.url', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find marker: spring.datasource.url" }
$targetPara = $searchRange.Paragraphs(1)
$targetPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="2B2B2B"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="CC7832"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>spring.datasource.url</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="808080"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>jdbc:postgresql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>://localhost:5432/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>postgres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="CC7832"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>spring.datasource.username</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="808080"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>postgres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="CC7832"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>spring.datasource</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="CC7832"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.password</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="808080"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>=</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="6A8759"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>password</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Item 5: 'Import collection from source : ...' -> split source/project with markers
$searchRange = $d.Content
$found = $searchRange.Find.Execute('Import collection from source', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find marker: Import collection from source" }
$targetPara = $searchRange.Paragraphs(1)
$targetPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">5. </w:t></w:r><w:r><w:t xml:space="preserve">Import collection from </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>source :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> /collection/</w:t></w:r><w:r><w:t xml:space="preserve">mini </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>project.postman_collection.json</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Remove the two empty paragraphs that followed item 5's paragraph
$searchRange2 = $d.Content
$found2 = $searchRange2.Find.Execute("Import collection from", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find item 5 paragraph for cleanup" }
$item5Para = $searchRange2.Paragraphs(1)
$emptyPara1 = $item5Para.Next()
$emptyPara1.Range.Delete()

$searchRange3 = $d.Content
$found3 = $searchRange3.Find.Execute("Import collection from", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$item5Para2 = $searchRange3.Paragraphs(1)
$emptyPara2 = $item5Para2.Next()
$emptyPara2.Range.Delete()
